$d = $word.ActiveDocument

# Locate the paragraph that contains the original sentence.
$para = $d.Paragraphs(1)
$full = $para.Range
$fullText = $full.Text

$target = "a branch"
$idx = $fullText.IndexOf($target)

if ($idx -ge 0) {
    $start = $full.Start

    # Range that exactly covers "a branch" inside the run.
    $rTarget = $d.Range($start + $idx, $start + $idx + $target.Length)

    # Replace "a branch" with "TESTING EDIT".
    $rTarget.Text = "TESTING EDIT"

    # Re-grab the range now occupied by the freshly inserted text.
    $newLen = "TESTING EDIT".Length
    $rNew = $d.Range($start + $idx, $start + $idx + $newLen)

    # Temporarily bookmarking the inserted text forces Word to keep it
    # in its own run (separate from the surrounding text) instead of
    # re-merging it into the neighbouring runs; deleting the bookmark
    # afterwards leaves the run boundaries (and therefore the three
    # separate <w:r> elements) intact without any left-over formatting
    # or bookmark markup.
    $bm = $d.Bookmarks.Add("__tmp_split_mark", $rNew)
    $d.Bookmarks("__tmp_split_mark").Delete()
}

Write-Output $d.Paragraphs(1).Range.Text
